# Kardex report update - version completa web febrero 2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: product description (merged B1:E1) ---
$ws.Range("B1").Value = "Estuche Spigen Liquid Air Apple iPhone XR - Negro"

# --- Row 2: CODIGO (barcode) / UBICACION ---
$ws.Range("B2").Value = 8809613763935
$ws.Range("E2").Value = "3C"

# --- Row 3: SKU / EXISTENCIA ---
$ws.Range("B3").Value = "064CS24872"
$ws.Range("E3").Value = 5

# --- Row 5: movement line ---
$ws.Range("A5").Value = 44221
$ws.Range("B5").Value = 0.44732638888889
$ws.Range("C5").Value = 4173011669

# --- Row 6: movement line (now the last data row, replacing former row 6) ---
$ws.Range("A6").Value = 44221
$ws.Range("B6").Value = 0.44269675925926
$ws.Range("C6").Value = "Cargue Inicial"
$ws.Range("D6").Value = 6
$ws.Range("E6").ClearContents()

# --- Remove former rows 7 and 8 (data no longer present) ---
$ws.Rows("7:8").Delete()

# --- Update selection to match the new last row ---
$ws.Range("A6").Select() | Out-Null
